# "bijna hele website responsive gemaakt"
# Duplicate the "Pakket Brons" package info (rows 1-6, column A) further down
# the sheet (rows 16-21) as plain stacked text, presumably for a responsive /
# mobile layout of the pricing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("A16").Value = "Pakket Brons"
$ws.Range("A17").Value = "10 GB HDD"
$ws.Range("A18").Value = "512 MB Ram"
$ws.Range("A19").Value = "1 Core"
$ws.Range("A20").Value = "100 GB Dataverkeer Per maand"
$ws.Range("A21").Value = "5,95 Per Maand"

$ws.Range("A1:G1").Select()
